$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 387, shifting existing rows 387:478 down to 388:479.
$ws.Rows("387:387").Insert()

# Populate the newly inserted row 387 with the new record.
$ws.Cells.Item(387, 1).Value2  = 3
$ws.Cells.Item(387, 2).Value2  = "Femacal de La Calera"
$ws.Cells.Item(387, 3).Value2  = "Coquimbo"
$ws.Cells.Item(387, 4).Value2  = 44722
$ws.Cells.Item(387, 5).Value2  = 5
$ws.Cells.Item(387, 6).Value2  = 100112021
$ws.Cells.Item(387, 7).Value2  = "Ají"
$ws.Cells.Item(387, 8).Value2  = "Inferno"
$ws.Cells.Item(387, 9).Value2  = "Primera"
$ws.Cells.Item(387, 10).Value2 = 78
$ws.Cells.Item(387, 11).Value2 = 21000
$ws.Cells.Item(387, 12).Value2 = 22000
$ws.Cells.Item(387, 13).Value2 = 21487
$ws.Cells.Item(387, 14).Value2 = "$/caja 15 kilos"
$ws.Cells.Item(387, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(387, 16).Value2 = 1432
$ws.Cells.Item(387, 17).Value2 = 15
$ws.Cells.Item(387, 18).Value2 = "Hortaliza"
